# Weekly update: a new week's record is inserted at row 3 (just below the
# most recent record in row 2). The existing records that occupied rows
# 3-23 each shift down by one row (to rows 4-24). Row 25 and below are
# untouched. Only the columns that actually carry record-specific data
# (D = Fecha, J = Volumen, K = Precio minimo, L = Precio maximo,
# M = Precio promedio ponderado, P = Precio $/Kg) need to move; the rest
# of each row (region, product, unit, origin, etc.) is identical for every
# row in this subset sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "P")

# Capture the current (pre-edit) values for rows 3..23 before we start
# overwriting anything, so the downward shift doesn't clobber data we
# still need to read.
$captured = @{}
for ($r = 3; $r -le 23; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $captured[$r] = $row
}

# Shift rows 3..23 down into rows 4..24.
for ($r = 23; $r -ge 3; $r--) {
    $target = $r + 1
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $captured[$r][$c]
    }
}

# New record for row 3 (the latest week).
$ws.Range("D3").Value = 44631
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 17000
$ws.Range("M3").Value = 16500
$ws.Range("P3").Value = 1269
